# Hindalco prices: add latest day's row (31-12-2025) at the top of the
# table, pushing every existing row down by one. The table is a running
# daily log ordered most-recent-first in row 2, so this is modeled as an
# insert of a new row 2 followed by populating it with the newest data
# point (date/price/circular taken from site's 25-12-2025 circular, same
# as the prior top row before the shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 2..203 down to 3..204, creating a fresh blank row 2.
$ws.Rows.Item(2).Insert()

# Populate the new top row with the latest entry.
$ws.Range("A2").Value = "31-12-2025"
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 312.5
$ws.Range("E2").Value = "25.12.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-25-december-2025.pdf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-25-december-2025.pdf")

# The inserted row inherited bold/centered header-style formatting from
# row 1 above it; restore the plain data-row look by copying formats
# from the row right below (now row 3, an untouched former data row).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
